# Apply the updated crypto price/volume figures scraped on
# Thu Apr 25 05:07:43 UTC 2024 (GitHub Actions refresh).
#
# The D (Price) / E (Volume(1h)) columns hold plain text in the
# source workbook (things like "64.416.17" or "  -3.83%  " are
# not valid numbers/percentages to Excel). Assigning through
# .Value with a leading apostrophe forces Excel's quote-prefix
# (text) interpretation so values such as "606.03" are kept as
# the literal string "606.03" instead of being silently coerced
# into a floating point number. Resetting .Style back to
# 'Normal' immediately afterwards clears the quote-prefix flag
# again so the cell format matches the original (unstyled) cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = '''64.304.28'
$c.Style = "Normal"
$c = $ws.Range('E2')
$c.Value = '''  -3.98%  '
$c.Style = "Normal"

$c = $ws.Range('D3')
$c.Value = '''3.156.70'
$c.Style = "Normal"
$c = $ws.Range('E3')
$c.Value = '''  -3.34%  '
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.Value = '''  +0.13%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.Value = '''606.03'
$c.Style = "Normal"
$c = $ws.Range('E5')
$c.Value = '''  -0.12%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.Value = '''147.20'
$c.Style = "Normal"
$c = $ws.Range('E6')
$c.Value = '''  -7.41%  '
$c.Style = "Normal"

$c = $ws.Range('E7')
$c.Value = '''  +0.12%  '
$c.Style = "Normal"

$c = $ws.Range('D8')
$c.Value = '''3.151.89'
$c.Style = "Normal"
$c = $ws.Range('E8')
$c.Value = '''  -3.55%  '
$c.Style = "Normal"

$c = $ws.Range('E9')
$c.Value = '''  -4.42%  '
$c.Style = "Normal"

$c = $ws.Range('E10')
$c.Value = '''  -7.61%  '
$c.Style = "Normal"

$c = $ws.Range('D11')
$c.Value = '''5.52'
$c.Style = "Normal"
$c = $ws.Range('E11')
$c.Value = '''  -6.94%  '
$c.Style = "Normal"

$c = $ws.Range('E12')
$c.Value = '''  -6.79%  '
$c.Style = "Normal"

$c = $ws.Range('D13')
$c.Value = '''0.0000250'
$c.Style = "Normal"
$c = $ws.Range('E13')
$c.Value = '''  -8.42%  '
$c.Style = "Normal"

$c = $ws.Range('D14')
$c.Value = '''36.13'
$c.Style = "Normal"
$c = $ws.Range('E14')
$c.Value = '''  -9.09%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.Value = '''3.675.08'
$c.Style = "Normal"
$c = $ws.Range('E15')
$c.Value = '''  -3.46%  '
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.Value = '''64.315.15'
$c.Style = "Normal"
$c = $ws.Range('E16')
$c.Value = '''  -3.98%  '
$c.Style = "Normal"

$c = $ws.Range('E17')
$c.Value = '''  +0.21%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.Value = '''3.154.56'
$c.Style = "Normal"
$c = $ws.Range('E18')
$c.Value = '''  -3.44%  '
$c.Style = "Normal"

$c = $ws.Range('D19')
$c.Value = '''6.94'
$c.Style = "Normal"
$c = $ws.Range('E19')
$c.Value = '''  -6.67%  '
$c.Style = "Normal"

$c = $ws.Range('D20')
$c.Value = '''479.53'
$c.Style = "Normal"
$c = $ws.Range('E20')
$c.Value = '''  -6.23%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.Value = '''14.70'
$c.Style = "Normal"
$c = $ws.Range('E21')
$c.Value = '''  -5.05%  '
$c.Style = "Normal"

$c = $ws.Range('E22')
$c.Value = '''  -6.81%  '
$c.Style = "Normal"

$c = $ws.Range('D23')
$c.Value = '''7.68'
$c.Style = "Normal"
$c = $ws.Range('E23')
$c.Value = '''  -5.88%  '
$c.Style = "Normal"

$c = $ws.Range('D24')
$c.Value = '''13.72'
$c.Style = "Normal"
$c = $ws.Range('E24')
$c.Value = '''  -8.06%  '
$c.Style = "Normal"

$c = $ws.Range('D25')
$c.Value = '''83.61'
$c.Style = "Normal"
$c = $ws.Range('E25')
$c.Value = '''  -3.20%  '
$c.Style = "Normal"

$c = $ws.Range('E26')
$c.Value = '''  -0.03%  '
$c.Style = "Normal"

$c = $ws.Range('E27')
$c.Value = '''  -5.31%  '
$c.Style = "Normal"

$c = $ws.Range('D28')
$c.Value = '''8.42'
$c.Style = "Normal"
$c = $ws.Range('E28')
$c.Value = '''  -8.41%  '
$c.Style = "Normal"

$c = $ws.Range('D29')
$c.Value = '''2.19'
$c.Style = "Normal"
$c = $ws.Range('E29')
$c.Value = '''  -8.37%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.Value = '''6.83'
$c.Style = "Normal"
$c = $ws.Range('E30')
$c.Value = '''  -2.08%  '
$c.Style = "Normal"

$c = $ws.Range('D31')
$c.Value = '''0.113'
$c.Style = "Normal"
$c = $ws.Range('E31')
$c.Value = '''  -35.40%  '
$c.Style = "Normal"

$c = $ws.Range('D32')
$c.Value = '''2.75'
$c.Style = "Normal"
$c = $ws.Range('E32')
$c.Value = '''  -5.36%  '
$c.Style = "Normal"

$c = $ws.Range('E33')
$c.Value = '''  +0.03%  '
$c.Style = "Normal"

$c = $ws.Range('D34')
$c.Value = '''26.19'
$c.Style = "Normal"
$c = $ws.Range('E34')
$c.Value = '''  -7.95%  '
$c.Style = "Normal"

$c = $ws.Range('E35')
$c.Value = '''  -5.21%  '
$c.Style = "Normal"

$c = $ws.Range('D36')
$c.Value = '''54.26'
$c.Style = "Normal"
$c = $ws.Range('E36')
$c.Value = '''  -2.55%  '
$c.Style = "Normal"

$c = $ws.Range('D37')
$c.Value = '''5.98'
$c.Style = "Normal"
$c = $ws.Range('E37')
$c.Value = '''  -7.50%  '
$c.Style = "Normal"

$c = $ws.Range('D38')
$c.Value = '''0.0₃0712'
$c.Style = "Normal"
$c = $ws.Range('E38')
$c.Value = '''  -11.32%  '
$c.Style = "Normal"

$c = $ws.Range('D39')
$c.Value = '''450.42'
$c.Style = "Normal"
$c = $ws.Range('E39')
$c.Value = '''  -9.62%  '
$c.Style = "Normal"

$c = $ws.Range('D40')
$c.Value = '''2.90'
$c.Style = "Normal"
$c = $ws.Range('E40')
$c.Value = '''  -13.89%  '
$c.Style = "Normal"

$c = $ws.Range('D41')
$c.Value = '''0.0395'
$c.Style = "Normal"
$c = $ws.Range('E41')
$c.Value = '''  -8.37%  '
$c.Style = "Normal"

$c = $ws.Range('E42')
$c.Value = '''  -8.39%  '
$c.Style = "Normal"

$c = $ws.Range('D43')
$c.Value = '''8.43'
$c.Style = "Normal"
$c = $ws.Range('E43')
$c.Value = '''  -4.74%  '
$c.Style = "Normal"

$c = $ws.Range('D44')
$c.Value = '''2.840.11'
$c.Style = "Normal"
$c = $ws.Range('E44')
$c.Value = '''  -4.02%  '
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.Value = '''0.267'
$c.Style = "Normal"
$c = $ws.Range('E45')
$c.Value = '''  -10.67%  '
$c.Style = "Normal"

$c = $ws.Range('E46')
$c.Value = '''  -9.46%  '
$c.Style = "Normal"

$c = $ws.Range('D47')
$c.Value = '''26.38'
$c.Style = "Normal"
$c = $ws.Range('E47')
$c.Value = '''  -8.55%  '
$c.Style = "Normal"

$c = $ws.Range('E48')
$c.Value = '''  -0.03%  '
$c.Style = "Normal"

$c = $ws.Range('E49')
$c.Value = '''  -5.33%  '
$c.Style = "Normal"

$c = $ws.Range('E50')
$c.Value = '''  -4.97%  '
$c.Style = "Normal"

$c = $ws.Range('D51')
$c.Value = '''118.23'
$c.Style = "Normal"
$c = $ws.Range('E51')
$c.Value = '''  -3.04%  '
$c.Style = "Normal"
